$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 890.9091
$ws.Range("I6").Value = 121.42857
$ws.Range("J6").Value = 5200
$ws.Range("K6").Value = 364.28571
$ws.Range("L6").Value = 15600
$ws.Range("M6").Value = -252.28571
$ws.Range("N6").Value = -15824
$ws.Range("H40").Value = 1500
$ws.Range("J40").Value = 1500
$ws.Range("L40").Value = 1500
$ws.Range("N40").Value = -1850
$ws.Range("H68").Value = 32800
$ws.Range("J68").Value = 32800
$ws.Range("L68").Value = 32800
$ws.Range("N68").Value = -34298
$ws.Range("H71").Value = 32800
$ws.Range("J71").Value = 32800
$ws.Range("L71").Value = 98400
$ws.Range("N71").Value = -105888
$ws.Range("H74").Value = 4215.4
$ws.Range("I74").Value = 3782.2856
$ws.Range("J74").Value = 5226
$ws.Range("K74").Value = 3782.2856
$ws.Range("L74").Value = 5226
$ws.Range("M74").Value = -2846.2856
$ws.Range("N74").Value = -7098
$ws.Range("H77").Value = 4215.4
$ws.Range("I77").Value = 3782.2856
$ws.Range("J77").Value = 5226
$ws.Range("K77").Value = 18911.428
$ws.Range("L77").Value = 26130
$ws.Range("M77").Value = -14231.428
$ws.Range("N77").Value = -35490
$ws.Range("H100").Value = 1516.3334
$ws.Range("I100").Value = 1410.3572
$ws.Range("J100").Value = 3000
$ws.Range("K100").Value = 1410.3572
$ws.Range("L100").Value = 3000
$ws.Range("M100").Value = -869.3571999999999
$ws.Range("N100").Value = -4082
$ws.Range("H116").Value = 2601
$ws.Range("I116").Value = 2376.25
$ws.Range("J116").Value = 3500
$ws.Range("K116").Value = 2376.25
$ws.Range("L116").Value = 3500
$ws.Range("M116").Value = 1065.75
$ws.Range("N116").Value = -10384
$ws.Range("H132").Value = 1443.6
$ws.Range("I132").Value = 1383.322
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 4149.965999999999
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -1619.965999999999
$ws.Range("N132").Value = -20060
$ws.Range("H138").Value = 1875.8163
$ws.Range("I138").Value = 1435.6
$ws.Range("J138").Value = 3832.3333
$ws.Range("K138").Value = 4306.799999999999
$ws.Range("L138").Value = 11496.9999
$ws.Range("M138").Value = 833.2000000000007
$ws.Range("N138").Value = -21776.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 115
$ws.Range("I4").Value = 30
$ws.Range("J4").Value = 200
$ws.Range("K4").Value = 30
$ws.Range("L4").Value = 200
$ws.Range("M4").Value = 86
$ws.Range("N4").Value = -432
$ws.Range("H14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("L14").ClearContents()
$ws.Range("N14").Value = 0
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").ClearContents()
$ws.Range("N15").Value = 0
$ws.Range("H21").Value = 34338.332
$ws.Range("I21").Value = 12007.5
$ws.Range("J21").Value = 79000
$ws.Range("K21").Value = 12007.5
$ws.Range("L21").Value = 79000
$ws.Range("M21").Value = -11633.5
$ws.Range("N21").Value = -79748
$ws.Range("H63").Value = 116827.625
$ws.Range("I63").Value = 302208
$ws.Range("J63").Value = 5599.4
$ws.Range("K63").Value = 302208
$ws.Range("L63").Value = 5599.4
$ws.Range("M63").Value = -301522
$ws.Range("N63").Value = -6971.4
$ws.Range("H66").Value = 116827.625
$ws.Range("I66").Value = 302208
$ws.Range("J66").Value = 5599.4
$ws.Range("K66").Value = 1511040
$ws.Range("L66").Value = 27997
$ws.Range("M66").Value = -1507608
$ws.Range("N66").Value = -34861
$ws.Range("H74").Value = 2309.1316
$ws.Range("I74").Value = 2020.3572
$ws.Range("K74").Value = 2020.3572
$ws.Range("M74").Value = -1146.3572
$ws.Range("H77").Value = 2309.1316
$ws.Range("I77").Value = 2020.3572
$ws.Range("K77").Value = 10101.786
$ws.Range("M77").Value = -5733.786
$ws.Range("H122").Value = 2900.8572
$ws.Range("I122").Value = 2237.4546
$ws.Range("K122").Value = 6712.3638
$ws.Range("M122").Value = -4262.3638
$ws.Range("H132").Value = 3739.878
$ws.Range("I132").Value = 2608.2593
$ws.Range("J132").Value = 5922.2856
$ws.Range("K132").Value = 7824.777900000001
$ws.Range("L132").Value = 17766.8568
$ws.Range("M132").Value = -5294.777900000001
$ws.Range("N132").Value = -22826.8568

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 5752
$ws.Range("I8").Value = 10004
$ws.Range("J8").Value = 1500
$ws.Range("K8").Value = 10004
$ws.Range("L8").Value = 1500
$ws.Range("M8").Value = -9864
$ws.Range("N8").Value = -1780
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("L10").ClearContents()
$ws.Range("M10").ClearContents()
$ws.Range("N10").Value = 0
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("M15").ClearContents()
$ws.Range("H35").Value = 33332.332
$ws.Range("J35").Value = 33332.332
$ws.Range("L35").Value = 33332.332
$ws.Range("N35").Value = -33952.332
$ws.Range("H82").Value = 30642.166
$ws.Range("I82").Value = 6928.5
$ws.Range("J82").Value = 42499
$ws.Range("K82").Value = 6928.5
$ws.Range("L82").Value = 42499
$ws.Range("M82").Value = -6545.5
$ws.Range("N82").Value = -43265
$ws.Range("H85").Value = 30642.166
$ws.Range("I85").Value = 6928.5
$ws.Range("J85").Value = 42499
$ws.Range("K85").Value = 6928.5
$ws.Range("L85").Value = 42499
$ws.Range("M85").Value = -5602.5
$ws.Range("N85").Value = -45151
$ws.Range("H134").Value = 2487.95
$ws.Range("I134").Value = 1850.2667
$ws.Range("J134").Value = 4401
$ws.Range("K134").Value = 5550.800099999999
$ws.Range("L134").Value = 13203
$ws.Range("M134").Value = -3015.800099999999
$ws.Range("N134").Value = -18273

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 9712.916999999999
$ws.Range("J3").Value = 9712.916999999999
$ws.Range("L3").Value = 9712.916999999999
$ws.Range("N3").Value = -9938.916999999999
$ws.Range("H31").Value = 7469.1396
$ws.Range("I31").Value = 1639.9231
$ws.Range("J31").Value = 9995.134
$ws.Range("K31").Value = 1639.9231
$ws.Range("L31").Value = 9995.134
$ws.Range("M31").Value = -1344.9231
$ws.Range("N31").Value = -10585.134
$ws.Range("H34").Value = 7469.1396
$ws.Range("I34").Value = 1639.9231
$ws.Range("J34").Value = 9995.134
$ws.Range("K34").Value = 1639.9231
$ws.Range("L34").Value = 9995.134
$ws.Range("M34").Value = -1437.9231
$ws.Range("N34").Value = -10399.134
$ws.Range("H99").Value = 1966.4348
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 1966.4348
$ws.Range("K99").Value = 0
$ws.Range("L99").ClearContents()
$ws.Range("M99").Value = 1966.4348
$ws.Range("N99").Value = -4962.4348
$ws.Range("H126").Value = 1966.4348
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 1966.4348
$ws.Range("K126").Value = 0
$ws.Range("L126").ClearContents()
$ws.Range("M126").Value = 5899.3044
$ws.Range("N126").Value = -10839.3044
$ws.Range("H132").Value = 21606130
$ws.Range("I132").Value = 33334542
$ws.Range("J132").Value = 6945614.5
$ws.Range("K132").Value = 100003626
$ws.Range("L132").Value = 20836843.5
$ws.Range("M132").Value = -100001096
$ws.Range("N132").Value = -20841903.5
$ws.Range("H134").Value = 4308.081
$ws.Range("I134").Value = 3925.1428
$ws.Range("J134").Value = 11009.5
$ws.Range("K134").Value = 11775.4284
$ws.Range("L134").Value = 33028.5
$ws.Range("M134").Value = -9240.428400000001
$ws.Range("N134").Value = -38098.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 700
$ws.Range("I87").Value = 700
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 2100
$ws.Range("L87").Value = 0
$ws.Range("M87").ClearContents()
$ws.Range("N87").Value = -852
$ws.Range("H90").Value = 700
$ws.Range("I90").Value = 700
$ws.Range("J90").Value = 0
$ws.Range("K90").Value = 6300
$ws.Range("L90").Value = 0
$ws.Range("M90").ClearContents()
$ws.Range("N90").Value = -60
$ws.Range("H131").Value = 1617.7693
$ws.Range("I131").Value = 550
$ws.Range("J131").Value = 1811.909
$ws.Range("K131").Value = 1650
$ws.Range("L131").Value = 5435.727000000001
$ws.Range("M131").Value = 3390
$ws.Range("N131").Value = -15515.727

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 1324901.1
$ws.Range("I11").Value = 2001499.9
$ws.Range("J11").Value = 310003
$ws.Range("K11").Value = 2001499.9
$ws.Range("L11").Value = 310003
$ws.Range("M11").Value = -2001360.9
$ws.Range("N11").Value = -310281
$ws.Range("H102").Value = 1644.6666
$ws.Range("I102").Value = 1424.5454
$ws.Range("K102").Value = 1424.5454
$ws.Range("M102").Value = 197.4546
$ws.Range("H132").Value = 3205.9546
$ws.Range("I132").Value = 2896.1353
$ws.Range("J132").Value = 4843.5713
$ws.Range("K132").Value = 8688.4059
$ws.Range("L132").Value = 14530.7139
$ws.Range("M132").Value = -6158.4059
$ws.Range("N132").Value = -19590.7139

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 103199.9
$ws.Range("I40").Value = 128000.5
$ws.Range("J40").Value = 3997.5
$ws.Range("K40").Value = 128000.5
$ws.Range("L40").Value = 3997.5
$ws.Range("M40").Value = -127864.5
$ws.Range("N40").Value = -4269.5
$ws.Range("H100").Value = 3860
$ws.Range("I100").Value = 4980
$ws.Range("J100").Value = 3300
$ws.Range("K100").Value = 4980
$ws.Range("L100").Value = 3300
$ws.Range("M100").Value = -4439
$ws.Range("N100").Value = -4382
$ws.Range("H132").Value = 2478.4546
$ws.Range("I132").Value = 2063.9644
$ws.Range("J132").Value = 4799.6
$ws.Range("K132").Value = 6191.8932
$ws.Range("L132").Value = 14398.8
$ws.Range("M132").Value = -3661.8932
$ws.Range("N132").Value = -19458.8
$ws.Range("H136").Value = 3473849
$ws.Range("I136").Value = 1242.2609
$ws.Range("J136").Value = 6668647.5
$ws.Range("K136").Value = 3726.7827
$ws.Range("L136").Value = 20005942.5
$ws.Range("M136").Value = -1176.7827
$ws.Range("N136").Value = -20011042.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 3341.5454
$ws.Range("I96").Value = 2053.8572
$ws.Range("J96").Value = 5595
$ws.Range("K96").Value = 2053.8572
$ws.Range("L96").Value = 5595
$ws.Range("M96").Value = -680.8571999999999
$ws.Range("N96").Value = -8341
$ws.Range("H136").Value = 2411.302
$ws.Range("I136").Value = 1972.4147
$ws.Range("K136").Value = 5917.2441
$ws.Range("M136").Value = -3367.2441
